$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '28.637.06'; E = '  +1.67%  '; DNumeric = $false },
    @{ Row = 3; D = '1.868.11'; E = '  +1.64%  '; DNumeric = $false },
    @{ Row = 4; D = '1.006'; E = '  +0.36%  '; DNumeric = $true },
    @{ Row = 5; D = '324.04'; E = '  -0.71%  '; DNumeric = $true },
    @{ Row = 6; D = '1.006'; E = '  +0.41%  '; DNumeric = $true },
    @{ Row = 7; D = '0.4549'; E = '  -1.87%  '; DNumeric = $true },
    @{ Row = 8; D = '0.3831'; E = '  -0.99%  '; DNumeric = $true },
    @{ Row = 9; D = '0.07790'; E = '  -0.84%  '; DNumeric = $true },
    @{ Row = 10; D = '0.9845'; E = '  +2.22%  '; DNumeric = $true },
    @{ Row = 11; D = '21.48'; E = '  -2.56%  '; DNumeric = $true },
    @{ Row = 12; D = '1.892.32'; E = '  +1.19%  '; DNumeric = $false },
    @{ Row = 13; D = '6.916'; E = '  +0.36%  '; DNumeric = $true },
    @{ Row = 14; D = '5.631'; E = '  -0.98%  '; DNumeric = $true },
    @{ Row = 15; D = '0.06953'; E = '  +1.22%  '; DNumeric = $true },
    @{ Row = 16; D = '87.94'; E = '  -0.62%  '; DNumeric = $true },
    @{ Row = 17; D = '1.007'; E = '  +0.42%  '; DNumeric = $true },
    @{ Row = 18; D = '0.000009931'; E = '  -0.07%  '; DNumeric = $true },
    @{ Row = 19; D = '16.65'; E = '  -0.42%  '; DNumeric = $true },
    @{ Row = 20; D = '1.005'; E = '  +0.27%  '; DNumeric = $true },
    @{ Row = 21; D = '28.660.77'; E = '  +1.70%  '; DNumeric = $false },
    @{ Row = 22; D = '5.236'; E = '  -1.24%  '; DNumeric = $true },
    @{ Row = 23; D = '10.89'; E = '  -1.26%  '; DNumeric = $true },
    @{ Row = 24; D = '2.110'; E = '  +0.75%  '; DNumeric = $true },
    @{ Row = 25; D = '2.110.77'; E = '  +1.51%  '; DNumeric = $false },
    @{ Row = 26; D = '152.75'; E = '  -1.06%  '; DNumeric = $true },
    @{ Row = 27; D = '19.12'; E = '  -0.28%  '; DNumeric = $true },
    @{ Row = 28; D = '5.684'; E = '  -0.63%  '; DNumeric = $true },
    @{ Row = 29; D = '1.931'; E = '  -2.06%  '; DNumeric = $true },
    @{ Row = 30; D = '118.04'; E = '  -0.98%  '; DNumeric = $true },
    @{ Row = 31; D = '0.09265'; E = '  -0.01%  '; DNumeric = $true },
    @{ Row = 32; D = '0.9049'; E = '  -3.40%  '; DNumeric = $true },
    @{ Row = 33; D = '5.271'; E = '  -0.22%  '; DNumeric = $true },
    @{ Row = 34; D = '1.315'; E = '  -0.55%  '; DNumeric = $true },
    @{ Row = 35; D = '3.298'; E = '  -0.76%  '; DNumeric = $true },
    @{ Row = 36; D = '0.05705'; E = '  -2.24%  '; DNumeric = $true },
    @{ Row = 37; D = '1.138'; E = '  -0.08%  '; DNumeric = $true },
    @{ Row = 38; D = '0.02062'; E = '  -3.02%  '; DNumeric = $true },
    @{ Row = 39; D = '7.654'; E = '  -1.43%  '; DNumeric = $true },
    @{ Row = 40; D = '0.5570'; E = '  -0.49%  '; DNumeric = $true },
    @{ Row = 41; D = '0.1770'; E = '  +0.58%  '; DNumeric = $true },
    @{ Row = 42; D = '9.620'; E = '  -2.98%  '; DNumeric = $true },
    @{ Row = 43; D = '0.07085'; E = '  -3.41%  '; DNumeric = $true },
    @{ Row = 44; D = '11.58'; E = '  -0.70%  '; DNumeric = $true },
    @{ Row = 45; D = '0.5229'; E = '  -0.85%  '; DNumeric = $true },
    @{ Row = 46; D = '2.124'; E = '  -0.73%  '; DNumeric = $true },
    @{ Row = 47; D = '1.816'; E = '  -1.37%  '; DNumeric = $true },
    @{ Row = 48; D = '112.37'; E = '  -1.46%  '; DNumeric = $true },
    @{ Row = 49; D = '1.099'; E = '  -3.70%  '; DNumeric = $true },
    @{ Row = 50; D = '2.424'; E = '  +4.29%  '; DNumeric = $true },
    @{ Row = 51; D = '1.005'; E = '  +0.38%  '; DNumeric = $true }
)

foreach ($u in $updates) {
    $dCell = $ws.Cells.Item($u.Row, 4)
    if ($u.DNumeric) {
        # Force text so numeric-looking strings (e.g. "1.006") are not
        # auto-converted to numbers by Excel, matching the original
        # text representation of this column.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    } else {
        $dCell.Value = $u.D
    }

    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
